$d = $word.ActiveDocument

$wrapperPre = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>'
$wrapperPost = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$bodies = @(
  '<w:p><w:r><w:t xml:space="preserve">This change is made in </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>dev</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> branch.</w:t></w:r></w:p>',
  '<w:p><w:r><w:t xml:space="preserve">  </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t xml:space="preserve">Pushing on remote repo using </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>dev</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> branch.</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p>',
  '<w:p><w:proofErr w:type="spellStart"/><w:r><w:t>Aaaaaaaaaaaaaaaaaaaaaaaaaaaa</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>',
  '<w:p><w:proofErr w:type="spellStart"/><w:r><w:t>Bbbbbbbbbbbbbbbbbbbbbbbbbbbbb</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>',
  '<w:p><w:proofErr w:type="spellStart"/><w:r><w:t>Ccccccccccccccccccccccccccccc</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>',
  '<w:p><w:proofErr w:type="spellStart"/><w:r><w:t>S</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>;dfjklsdjfkjk</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>Modified second file to see the functionality.</w:t></w:r></w:p>',
  '<w:p><w:r><w:t xml:space="preserve">This change is made in </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>dev</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> branch.</w:t></w:r></w:p>',
  '<w:p><w:r><w:t xml:space="preserve">  </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t xml:space="preserve">Pushing on remote repo using </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>dev</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> branch.</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p>',
  '<w:p><w:proofErr w:type="spellStart"/><w:r><w:t>Aaaaaaaaaaaaaaaaaaaaaaaaaaaa</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>',
  '<w:p><w:proofErr w:type="spellStart"/><w:r><w:t>Bbbbbbbbbbbbbbbbbbbbbbbbbbbbb</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>',
  '<w:p><w:proofErr w:type="spellStart"/><w:r><w:t>Ccccccccccccccccccccccccccccc</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>',
  '<w:p><w:proofErr w:type="spellStart"/><w:r><w:t>S</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>;dfjklsdjfkjk</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/></w:p>',
  '<w:p><w:r><w:t>Modified second file to see the functionality.</w:t></w:r></w:p>',
  '<w:p><w:r><w:t xml:space="preserve">This change is made in </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>dev</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> branch.</w:t></w:r></w:p>',
  '<w:p><w:r><w:t xml:space="preserve">  </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t xml:space="preserve">Pushing on remote repo using </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>dev</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> branch.</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p>',
  '<w:p><w:proofErr w:type="spellStart"/><w:r><w:t>Aaaaaaaaaaaaaaaaaaaaaaaaaaaa</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>',
  '<w:p><w:proofErr w:type="spellStart"/><w:r><w:t>Bbbbbbbbbbbbbbbbbbbbbbbbbbbbb</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>',
  '<w:p><w:proofErr w:type="spellStart"/><w:r><w:t>Ccccccccccccccccccccccccccccc</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>',
  '<w:p><w:proofErr w:type="spellStart"/><w:r><w:t>S</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>;dfjklsdjfkjk</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/></w:p>',
  '<w:p><w:r><w:t>Modified second file to see the functionality.</w:t></w:r></w:p>',
  '<w:p><w:r><w:t xml:space="preserve">This change is made in </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>dev</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> branch.</w:t></w:r></w:p>',
  '<w:p><w:r><w:t xml:space="preserve">  </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t xml:space="preserve">Pushing on remote repo using </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>dev</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> branch.</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p>',
  '<w:p><w:proofErr w:type="spellStart"/><w:r><w:t>Aaaaaaaaaaaaaaaaaaaaaaaaaaaa</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>',
  '<w:p><w:proofErr w:type="spellStart"/><w:r><w:t>Bbbbbbbbbbbbbbbbbbbbbbbbbbbbb</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>',
  '<w:p><w:proofErr w:type="spellStart"/><w:r><w:t>Ccccccccccccccccccccccccccccc</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>',
  '<w:p><w:proofErr w:type="spellStart"/><w:r><w:t>S</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>;dfjklsdjfkjk</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/></w:p>',
  '<w:p><w:r><w:t>Modified second file to see the functionality.</w:t></w:r></w:p>',
  '<w:p><w:r><w:t xml:space="preserve">This change is made in </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>dev</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> branch.</w:t></w:r></w:p>',
  '<w:p><w:r><w:t xml:space="preserve">  </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t xml:space="preserve">Pushing on remote repo using </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>dev</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> branch.</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p>',
  '<w:p><w:proofErr w:type="spellStart"/><w:r><w:t>Aaaaaaaaaaaaaaaaaaaaaaaaaaaa</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>',
  '<w:p><w:proofErr w:type="spellStart"/><w:r><w:t>Bbbbbbbbbbbbbbbbbbbbbbbbbbbbb</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>',
  '<w:p><w:proofErr w:type="spellStart"/><w:r><w:t>Ccccccccccccccccccccccccccccc</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>',
  '<w:p><w:proofErr w:type="spellStart"/><w:r><w:t>S</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>;dfjklsdjfkjk</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/></w:p>',
  '<w:p><w:r><w:t>Modified second file to see the functionality.</w:t></w:r></w:p>',
  '<w:p><w:r><w:t xml:space="preserve">This change is made in </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>dev</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> branch.</w:t></w:r></w:p>',
  '<w:p><w:r><w:t xml:space="preserve">  </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t xml:space="preserve">Pushing on remote repo using </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>dev</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> branch.</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p>',
  '<w:p><w:proofErr w:type="spellStart"/><w:r><w:t>Aaaaaaaaaaaaaaaaaaaaaaaaaaaa</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>',
  '<w:p><w:proofErr w:type="spellStart"/><w:r><w:t>Bbbbbbbbbbbbbbbbbbbbbbbbbbbbb</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>',
  '<w:p><w:proofErr w:type="spellStart"/><w:r><w:t>Ccccccccccccccccccccccccccccc</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>',
  '<w:p><w:proofErr w:type="spellStart"/><w:r><w:t>S</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>;dfjklsdjfkjk</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/></w:p>'
)

# Step 1: replace the "This change is made in dev branch." paragraph (paragraph 2)
# in place with the first reworked paragraph (adds proofErr spans around "dev").
$para2 = $d.Paragraphs(2)
$r2 = $para2.Range
$r2.InsertXML($wrapperPre + $bodies[0] + $wrapperPost)

# Step 2: insert all the remaining new paragraphs, one at a time, immediately
# before the trailing bookmark paragraph (which always ends up as the last
# paragraph in the document after each insertion).
for ($i = 1; $i -lt $bodies.Length; $i++) {
    $lastIdx = $d.Paragraphs.Count
    $target = $d.Paragraphs($lastIdx)
    $rt = $target.Range
    $rt.Collapse(1)
    $rt.InsertXML($wrapperPre + $bodies[$i] + $wrapperPost)
}
